{"js": "// Commit: \"Added the acc Activation using email\"\n// Net content change: a new, empty paragraph is inserted at the very end of\n// the document body -- right after the \"Backend APIs should be available\n// for frontend integration\" line and right before the (pre-existing)\n// trailing empty paragraph. The new paragraph copies that trailing\n// paragraph's formatting: spacing-after = 0pt and a 12pt (sz/szCs 24) run\n// font size, with no list numbering and no explicit language tag.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst anchorText = \"Backend APIs should be available for frontend integration\";\n\n// Default: the very last paragraph of the document (the existing trailing\n// empty paragraph) -- insert the new paragraph immediately before it.\nlet reference = items[items.length - 1];\n\n// Prefer anchoring right after the known sentence if we can find it, which\n// is robust even if the document already picked up extra trailing\n// paragraphs elsewhere.\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text.trim() === anchorText) {\n    reference = items[i + 1];\n    break;\n  }\n}\n\nconst inserted = reference.insertParagraph(\"\", \"Before\");\ninserted.spaceAfter = 0;\ninserted.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Commit: \"Added the acc Activation using email\"\n# Net content change: a new, empty paragraph is inserted at the very end of\n# the document body -- right after the \"Backend APIs should be available\n# for frontend integration\" line and right before the (pre-existing)\n# trailing empty paragraph. The new paragraph copies that trailing\n# paragraph's formatting: spacing-after = 0pt and a 12pt run font size,\n# with no list numbering and no explicit language tag.\n\n$d = $word.ActiveDocument\n$anchorText = \"Backend APIs should be available for frontend integration\"\n\n$count = $d.Paragraphs.Count\n$referenceIndex = $count\nfor ($i = 1; $i -le $count - 1; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text.TrimEnd(\"`r\") -eq $anchorText) {\n    $referenceIndex = $i + 1\n    break\n  }\n}\n\n$reference = $d.Paragraphs.Item($referenceIndex)\n$r = $reference.Range\n$r.Collapse(1)   # wdCollapseStart\n$r.InsertParagraphBefore()\n\n# The freshly inserted paragraph now sits immediately before $reference;\n# re-fetch it by index and pin its formatting to match the sibling empty\n# paragraph that follows it (spacing after = 0pt, 12pt font).\n$insertedIndex = $referenceIndex\n$inserted = $d.Paragraphs.Item($insertedIndex)\n$inserted.SpaceAfter = 0\n$inserted.Range.Font.Size = 12\n"}
